$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 (ORM No 17041805 - RAM_MPP.doc reviewed by Kaung Myat Bo)
$ws.Range("B6").Value = 17041805
$ws.Range("C6").Value = "RAM_MPP.doc"
$ws.Range("D6").Value = "Kaung Myat Bo"
$ws.Range("E6").Value = 43207

# Row 7 (ORM No 17041806 - RAM_MWBS.xlsx reviewed by Kaung Myat Bo)
$ws.Range("B7").Value = 17041806
$ws.Range("C7").Value = "RAM_MWBS.xlsx"
$ws.Range("D7").Value = "Kaung Myat Bo"
$ws.Range("E7").Value = 43207

# Copy date formatting from an existing date cell so the new date cells
# reuse the same style (numFmtId 14) instead of creating a new style.
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E6:E7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the saved view state.
$ws.Range("L8").Select() | Out-Null

Write-Host "Edit applied"
